$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Range("H9").Value = 201.3077
$ws.Range("I9").Value = 27.25
$ws.Range("K9").Value = 27.25
$ws.Range("M9").Value = 141.75

$ws.Range("H17").Value = 2007.2222
$ws.Range("J17").Value = 2558.4546
$ws.Range("L17").Value = 7675.3638
$ws.Range("N17").Value = -8011.3638

$ws.Range("H40").Value = 2459.7856
$ws.Range("I40").Value = 1617.25
$ws.Range("K40").Value = 1617.25
$ws.Range("M40").Value = -1442.25

$ws.Range("H53").Value = 601.0833
$ws.Range("I53").Value = 225
$ws.Range("J53").Value = 635.2727
$ws.Range("K53").Value = 225
$ws.Range("L53").Value = 635.2727
$ws.Range("M53").Value = 412
$ws.Range("N53").Value = -1909.2727

$ws.Range("H88").Value = 6594.9
$ws.Range("I88").Value = 4460
$ws.Range("J88").Value = 7306.533
$ws.Range("K88").Value = 4460
$ws.Range("L88").Value = 7306.533
$ws.Range("M88").Value = -4054
$ws.Range("N88").Value = -8118.533

$ws.Range("H91").Value = 6594.9
$ws.Range("I91").Value = 4460
$ws.Range("J91").Value = 7306.533
$ws.Range("K91").Value = 4460
$ws.Range("L91").Value = 7306.533
$ws.Range("M91").Value = -3056
$ws.Range("N91").Value = -10114.533

$ws.Range("H113").Value = 4010.5
$ws.Range("I113").Value = 3955.625
$ws.Range("J113").Value = 4449.5
$ws.Range("K113").Value = 3955.625
$ws.Range("L113").Value = 4449.5
$ws.Range("M113").Value = -701.625
$ws.Range("N113").Value = -10957.5

$ws.Range("H116").Value = 11124.625
$ws.Range("I116").Value = 10499.5
$ws.Range("J116").Value = 11749.75
$ws.Range("K116").Value = 10499.5
$ws.Range("L116").Value = 11749.75
$ws.Range("M116").Value = -7057.5
$ws.Range("N116").Value = -18633.75

$ws.Range("H132").Value = 11115109
$ws.Range("I132").Value = 4496.857
$ws.Range("J132").Value = 50002250
$ws.Range("K132").Value = 13490.571
$ws.Range("L132").Value = 150006750
$ws.Range("M132").Value = -10960.571
$ws.Range("N132").Value = -150011810

$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Range("H45").Value = 3106.3914
$ws.Range("I45").Value = 3183.625
$ws.Range("J45").Value = 2929.8572
$ws.Range("K45").Value = 3183.625
$ws.Range("L45").Value = 2929.8572
$ws.Range("M45").Value = -2806.625
$ws.Range("N45").Value = -3683.8572

$ws.Range("H61").Value = 3127462.2
$ws.Range("I61").Value = 2850
$ws.Range("J61").Value = 4168999.8
$ws.Range("K61").Value = 2850
$ws.Range("L61").Value = 4168999.8
$ws.Range("M61").Value = -2638
$ws.Range("N61").Value = -4169423.8

$ws.Range("H124").Value = 36806.668
$ws.Range("J124").Value = 45210
$ws.Range("L124").Value = 45210
$ws.Range("N124").Value = -55030

$ws.Range("H127").Value = 80000
$ws.Range("J127").Value = 80000
$ws.Range("L127").Value = 80000
$ws.Range("N127").Value = -89920

$ws.Range("H132").Value = 2314.0908
$ws.Range("I132").Value = 2314.0908
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6942.2724
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = -4412.2724
$ws.Range("M132").ClearContents()

$ws.Range("H136").Value = 3127462.2
$ws.Range("I136").Value = 2850
$ws.Range("J136").Value = 4168999.8
$ws.Range("K136").Value = 8550
$ws.Range("L136").Value = 12506999.4
$ws.Range("M136").Value = -6000
$ws.Range("N136").Value = -12512099.4

$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Range("H86").Value = 16129887
$ws.Range("I86").Value = 630.9545000000001
$ws.Range("K86").Value = 630.9545000000001
$ws.Range("M86").Value = 492.0454999999999

$ws.Range("H89").Value = 16129887
$ws.Range("I89").Value = 630.9545000000001
$ws.Range("K89").Value = 3154.7725
$ws.Range("M89").Value = 2461.2275

$ws.Range("H99").Value = 4457.353
$ws.Range("I99").Value = 2447.1667
$ws.Range("K99").Value = 2447.1667
$ws.Range("M99").Value = -949.1667000000002

$ws.Range("H134").Value = 6251750
$ws.Range("I134").Value = 2333.3333
$ws.Range("K134").Value = 6999.999899999999
$ws.Range("M134").Value = -4464.999899999999

$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Range("H22").Value = 1507.7273
$ws.Range("I22").Value = 941.2857
$ws.Range("K22").Value = 941.2857
$ws.Range("M22").Value = -591.2857

$ws.Range("H132").Value = 7272.1113
$ws.Range("I132").Value = 2362.5
$ws.Range("K132").Value = 7087.5
$ws.Range("M132").Value = -4557.5

$ws.Range("H134").Value = 3737.4
$ws.Range("I134").Value = 3796.75
$ws.Range("J134").Value = 3500
$ws.Range("K134").Value = 11390.25
$ws.Range("L134").Value = 10500
$ws.Range("M134").Value = -8855.25
$ws.Range("N134").Value = -15570

$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Range("H33").Value = 280.84616
$ws.Range("J33").Value = 299.3
$ws.Range("L33").Value = 1795.8
$ws.Range("N33").Value = -2361.8

$ws.Range("H121").Value = 1719.7778
$ws.Range("I121").Value = 283
$ws.Range("K121").Value = 849
$ws.Range("M121").Value = 461

$ws.Range("H131").Value = 4766323.5
$ws.Range("I131").Value = 1450.5
$ws.Range("K131").Value = 4351.5
$ws.Range("M131").Value = 688.5

$ws.Range("H132").Value = 1766.5428
$ws.Range("J132").Value = 1977.826
$ws.Range("L132").Value = 17800.434
$ws.Range("N132").Value = -22860.434

$ws.Range("H140").Value = 2057.7144
$ws.Range("I140").Value = 2057.7144
$ws.Range("K140").Value = 6173.1432
$ws.Range("M140").Value = -993.1431999999995

$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Range("H80").Value = 66671332
$ws.Range("I80").Value = 125002250
$ws.Range("K80").Value = 125002250
$ws.Range("M80").Value = -125001252

$ws.Range("H83").Value = 66671332
$ws.Range("I83").Value = 125002250
$ws.Range("K83").Value = 625011250
$ws.Range("M83").Value = -625006258

$ws.Range("H122").Value = 3139.7334
$ws.Range("J122").Value = 6399.75
$ws.Range("L122").Value = 19199.25
$ws.Range("N122").Value = -24099.25

$ws.Range("H126").Value = 2160.762
$ws.Range("I126").Value = 1713
$ws.Range("K126").Value = 5139
$ws.Range("M126").Value = -2669

$ws.Range("H132").Value = 41669056
$ws.Range("I132").Value = 55557390
$ws.Range("J132").Value = 4054.5
$ws.Range("K132").Value = 166672170
$ws.Range("L132").Value = 12163.5
$ws.Range("M132").Value = -166669640
$ws.Range("N132").Value = -17223.5

$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Range("H7").Value = 3721.5186
$ws.Range("I7").Value = 3615.5715
$ws.Range("K7").Value = 3615.5715
$ws.Range("M7").Value = -3503.5715

$ws.Range("H100").Value = 52635584
$ws.Range("J100").Value = 33337910
$ws.Range("L100").Value = 33337910
$ws.Range("N100").Value = -33338992

$ws.Range("H126").Value = 3721.5186
$ws.Range("I126").Value = 3615.5715
$ws.Range("K126").Value = 10846.7145
$ws.Range("M126").Value = -8376.7145

$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Range("H2").Value = 3450897.8
$ws.Range("I2").Value = 6901334
$ws.Range("K2").Value = 6901334
$ws.Range("M2").Value = -6901222

$ws.Range("H23").Value = 234.14285
$ws.Range("I23").Value = 256.5
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = 256.5
$ws.Range("L23").Value = 100
$ws.Range("M23").Value = -27.5
$ws.Range("N23").Value = -558

$ws.Range("H45").Value = 19993.25
$ws.Range("J45").Value = 19993.25
$ws.Range("L45").Value = 19993.25
$ws.Range("N45").Value = -20975.25

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H132").Value = 18521732
$ws.Range("I132").Value = 25644246
$ws.Range("J132").Value = 3196.4
$ws.Range("K132").Value = 76932738
$ws.Range("L132").Value = 9589.200000000001
$ws.Range("M132").Value = -76930208
$ws.Range("N132").Value = -14649.2
